$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 42.714287
$ws.Range("I11").Value = 42.714287
$ws.Range("K11").Value = 42.714287
$ws.Range("M11").Value = 97.285713

$ws.Range("H12").Value = 723.7778
$ws.Range("I12").Value = 986.2727
$ws.Range("J12").Value = 311.2857
$ws.Range("K12").Value = 986.2727
$ws.Range("L12").Value = 311.2857
$ws.Range("M12").Value = -816.2727
$ws.Range("N12").Value = -651.2857

$ws.Range("H33").Value = 1151361.4
$ws.Range("I33").Value = 1438068.4
$ws.Range("K33").Value = 1438068.4
$ws.Range("M33").Value = -1437839.4

$ws.Range("H76").Value = 4970.074
$ws.Range("I76").Value = 4989.6
$ws.Range("K76").Value = 4989.6
$ws.Range("M76").Value = -4674.6

$ws.Range("H79").Value = 4970.074
$ws.Range("I79").Value = 4989.6
$ws.Range("K79").Value = 4989.6
$ws.Range("M79").Value = -3897.6

$ws.Range("H101").Value = 625524.7
$ws.Range("I101").Value = 769763.6
$ws.Range("K101").Value = 2309290.8
$ws.Range("M101").Value = -2307668.8

$ws.Range("H137").Value = 69883.13
$ws.Range("I137").Value = 2691
$ws.Range("J137").Value = 114677.89
$ws.Range("K137").Value = 8073
$ws.Range("L137").Value = 344033.67
$ws.Range("M137").Value = -5523
$ws.Range("N137").Value = -349133.67

$ws.Range("H138").Value = 3805.818
$ws.Range("I138").Value = 2421
$ws.Range("J138").Value = 4826.2104
$ws.Range("K138").Value = 7263
$ws.Range("L138").Value = 14478.6312
$ws.Range("M138").Value = -2123
$ws.Range("N138").Value = -24758.6312

$ws.Range("H141").Value = 41528.695
$ws.Range("I141").Value = 45288.332
$ws.Range("K141").Value = 135864.996
$ws.Range("M141").Value = -130684.996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 130955.41
$ws.Range("I32").Value = 121192.96
$ws.Range("J32").Value = 252986
$ws.Range("K32").Value = 121192.96
$ws.Range("L32").Value = 252986
$ws.Range("M32").Value = -120905.96
$ws.Range("N32").Value = -253560

$ws.Range("H122").Value = 16602.076
$ws.Range("I122").Value = 16602.076
$ws.Range("K122").Value = 49806.228
$ws.Range("M122").Value = -47356.228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2212.818
$ws.Range("I107").Value = 1225.5
$ws.Range("J107").Value = 3397.6
$ws.Range("K107").Value = 1225.5
$ws.Range("L107").Value = 3397.6
$ws.Range("M107").Value = 694.5
$ws.Range("N107").Value = -7237.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 535.1818
$ws.Range("I7").Value = 561.25
$ws.Range("K7").Value = 561.25
$ws.Range("M7").Value = -448.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 387.5
$ws.Range("I35").Value = 387.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1162.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -874.5
$ws.Range("N35").Value = $null

$ws.Range("H139").Value = 2187
$ws.Range("I139").Value = 2177.5
$ws.Range("J139").Value = 2225
$ws.Range("K139").Value = 6532.5
$ws.Range("L139").Value = 6675
$ws.Range("M139").Value = -1392.5
$ws.Range("N139").Value = -16955

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 3882.8333
$ws.Range("I36").Value = 6558.5
$ws.Range("J36").Value = 2545
$ws.Range("K36").Value = 6558.5
$ws.Range("L36").Value = 2545
$ws.Range("M36").Value = -6073.5
$ws.Range("N36").Value = -3515

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null

$ws.Range("H80").Value = 11270.357
$ws.Range("I80").Value = 15375.375
$ws.Range("J80").Value = 5797
$ws.Range("K80").Value = 15375.375
$ws.Range("L80").Value = 5797
$ws.Range("M80").Value = -14377.375
$ws.Range("N80").Value = -7793

$ws.Range("H83").Value = 11270.357
$ws.Range("I83").Value = 15375.375
$ws.Range("J83").Value = 5797
$ws.Range("K83").Value = 76876.875
$ws.Range("L83").Value = 28985
$ws.Range("M83").Value = -71884.875
$ws.Range("N83").Value = -38969

$ws.Range("H93").Value = 30251
$ws.Range("J93").Value = 30251
$ws.Range("L93").Value = 30251
$ws.Range("N93").Value = -33995

$ws.Range("H97").Value = 23925.822
$ws.Range("I97").Value = 35011.633
$ws.Range("K97").Value = 35011.633
$ws.Range("M97").Value = -34515.633

$ws.Range("H102").Value = 4167.8945
$ws.Range("I102").Value = 3580.7334
$ws.Range("K102").Value = 3580.7334
$ws.Range("M102").Value = -1958.7334

$ws.Range("H126").Value = 3040.2
$ws.Range("I126").Value = 2848.2856
$ws.Range("K126").Value = 8544.856800000001
$ws.Range("M126").Value = -6074.856800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6331.1665
$ws.Range("J16").Value = 6498
$ws.Range("L16").Value = 6498
$ws.Range("N16").Value = -6838

$ws.Range("H122").Value = 10114.4375
$ws.Range("I122").Value = 10845.071
$ws.Range("K122").Value = 32535.213
$ws.Range("M122").Value = -30085.213

$ws.Range("H132").Value = 4788.45
$ws.Range("I132").Value = 4651.6665
$ws.Range("J132").Value = 5198.8
$ws.Range("K132").Value = 13954.9995
$ws.Range("L132").Value = 15596.4
$ws.Range("M132").Value = -11424.9995
$ws.Range("N132").Value = -20656.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 40000
$ws.Range("J4").Value = 40000
$ws.Range("L4").Value = 40000
$ws.Range("N4").Value = -40226

$ws.Range("H68").Value = 181757
$ws.Range("I68").Value = 20000
$ws.Range("J68").Value = 262635.5
$ws.Range("K68").Value = 20000
$ws.Range("L68").Value = 262635.5
$ws.Range("M68").Value = -19189
$ws.Range("N68").Value = -264257.5

$ws.Range("H70").Value = 18971.54
$ws.Range("J70").Value = 18971.54
$ws.Range("L70").Value = 18971.54
$ws.Range("N70").Value = -19601.54

$ws.Range("H71").Value = 181757
$ws.Range("I71").Value = 20000
$ws.Range("J71").Value = 262635.5
$ws.Range("K71").Value = 60000
$ws.Range("L71").Value = 787906.5
$ws.Range("M71").Value = -55944
$ws.Range("N71").Value = -796018.5

$ws.Range("H73").Value = 18971.54
$ws.Range("J73").Value = 18971.54
$ws.Range("L73").Value = 18971.54
$ws.Range("N73").Value = -21155.54
